$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.940.52'
$ws.Range('E2').Value = '  +0.87%  '

$ws.Range('D3').Value = '2.597.14'
$ws.Range('E3').Value = '  +0.44%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.48%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.594'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.24%  '

$ws.Range('E9').Value = '  +2.24%  '

$ws.Range('E10').Value = '  +2.13%  '

$ws.Range('E11').Value = '  +0.43%  '

$ws.Range('E12').Value = '  +1.60%  '

$ws.Range('D13').Value = '3.050.59'
$ws.Range('E13').Value = '  +0.43%  '

$ws.Range('D14').Value = '60.941.61'
$ws.Range('E14').Value = '  +0.88%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.14%  '

$ws.Range('E16').Value = '  +0.96%  '

$ws.Range('D17').Value = '2.592.27'
$ws.Range('E17').Value = '  +0.10%  '

$ws.Range('E18').Value = '  -1.09%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '353.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.47%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.59%  '

$ws.Range('E21').Value = '  +1.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.77%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.426'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.45%  '

$ws.Range('E25').Value = '  -0.25%  '

$ws.Range('D26').Value = '2.710.72'
$ws.Range('E26').Value = '  +0.42%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.09%  '

$ws.Range('D28').Value = '0.0₃0844'
$ws.Range('E28').Value = '  +0.09%  '

$ws.Range('E29').Value = '  +0.59%  '

$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('E31').Value = '  +11.81%  '

$ws.Range('E32').Value = '  +0.17%  '

$ws.Range('E33').Value = '  +3.49%  '

$ws.Range('E35').Value = '  +5.11%  '

$ws.Range('E36').Value = '  +9.98%  '

$ws.Range('E37').Value = '  +1.24%  '

$ws.Range('E38').Value = '  +2.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.80'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.19%  '

$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.66%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.848'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.08%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '286.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.07%  '

$ws.Range('E43').Value = '  +1.63%  '

$ws.Range('E44').Value = '  +1.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0560'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.76%  '

$ws.Range('E46').Value = '  +0.00%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.06%  '

$ws.Range('E48').Value = '  +0.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0238'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.01%  '

$ws.Range('E50').Value = '  +0.04%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.95%  '
